$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example Scenarios")

$ws.Range("B63").Value = "statements:`n    - if:`n        match: true`n        block:`n            - return: string(""hello"")"
$ws.Range("A63").Value = "Simple if statement example"

$ws.Range("A64").Value = "If statement with elseif for fallback conditions"
$ws.Range("B64").Value = "statements:`n    - context: example := ""hello2""`n    - if:`n        match: context.example == ""hello""`n        block:`n            - return: string(""output"")`n        elseifs:`n            - match: context.example == ""hello1""`n              block:`n                - return: string(""output1"")`n            - match: context.example == ""hello2""`n              block:`n                - return: string(""output2"")`n    - return: """""

$ws.Range("A63:B64").WrapText = $true

$ws.Rows.Item(63).RowHeight = 72
$ws.Rows.Item(64).RowHeight = 201.6

$ws.Range("B67").Select()
